$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Insert()

# Copy full formatting (and everything) from column E to column D using Copy/PasteSpecial formats
$srcRange = $ws.Range("E5:E102")
$dstRange = $ws.Range("D5:D102")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

Write-Host "done"
